# natmiOut/YoungD4/LR-pairs_lrc2p/Col18a1-Gpc4.xlsx
# "Natmi following Dr Hou advice": the ligand-receptor table is regenerated with
# a new "M1" sending-cluster group (old "M2" -> "M1", old "sCs" -> "M2"), and a
# fresh "sCs" sending-cluster block is appended at the bottom (rows 14-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Col18a1/Gpc4)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col18a1"
$ws.Range("C2").Value = "Gpc4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 47.25665799999999
$ws.Range("H2").Value = 141.769974
$ws.Range("I2").Value = 0.6282707309614213
$ws.Range("J2").Value = 0.6282707309614212
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 15.740393
$ws.Range("N2").Value = 47.221179
$ws.Range("O2").Value = 0.2847568403735705
$ws.Range("P2").Value = 0.2847568403735705
$ws.Range("Q2").Value = 743.8383687865938
$ws.Range("R2").Value = 6694.545319079345
$ws.Range("S2").Value = 0.1789043882477679
$ws.Range("T2").Value = 0.1789043882477678

# Row 3: ECs -> FAPs (Col18a1/Gpc4)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col18a1"
$ws.Range("C3").Value = "Gpc4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 47.25665799999999
$ws.Range("H3").Value = 141.769974
$ws.Range("I3").Value = 0.6282707309614213
$ws.Range("J3").Value = 0.6282707309614212
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 33.82224299999999
$ws.Range("N3").Value = 101.466729
$ws.Range("O3").Value = 0.6118725911752718
$ws.Range("P3").Value = 0.6118725911752717
$ws.Range("Q3").Value = 1598.326170243894
$ws.Range("R3").Value = 14384.93553219504
$ws.Range("S3").Value = 0.3844216401129469
$ws.Range("T3").Value = 0.3844216401129468

# Row 4: ECs -> sCs (Col18a1/Gpc4)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col18a1"
$ws.Range("C4").Value = "Gpc4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 47.25665799999999
$ws.Range("H4").Value = 141.769974
$ws.Range("I4").Value = 0.6282707309614213
$ws.Range("J4").Value = 0.6282707309614212
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.713974666666666
$ws.Range("N4").Value = 17.141924
$ws.Range("O4").Value = 0.1033705684511578
$ws.Range("P4").Value = 0.1033705684511578
$ws.Range("Q4").Value = 270.0233466433306
$ws.Range("R4").Value = 2430.210119789976
$ws.Range("S4").Value = 0.06494470260070656
$ws.Range("T4").Value = 0.06494470260070655

# Row 5: FAPs -> ECs (Col18a1/Gpc4)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col18a1"
$ws.Range("C5").Value = "Gpc4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.28572166666667
$ws.Range("H5").Value = 33.85716499999999
$ws.Range("I5").Value = 0.1500421083721963
$ws.Range("J5").Value = 0.1500421083721963
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 15.740393
$ws.Range("N5").Value = 47.221179
$ws.Range("O5").Value = 0.2847568403735705
$ws.Range("P5").Value = 0.2847568403735705
$ws.Range("Q5").Value = 177.6416943219483
$ws.Range("R5").Value = 1598.775248897535
$ws.Range("S5").Value = 0.04272551670305546
$ws.Range("T5").Value = 0.04272551670305545

# Row 6: FAPs -> FAPs (Col18a1/Gpc4)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col18a1"
$ws.Range("C6").Value = "Gpc4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.28572166666667
$ws.Range("H6").Value = 33.85716499999999
$ws.Range("I6").Value = 0.1500421083721963
$ws.Range("J6").Value = 0.1500421083721963
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 33.82224299999999
$ws.Range("N6").Value = 101.466729
$ws.Range("O6").Value = 0.6118725911752718
$ws.Range("P6").Value = 0.6118725911752717
$ws.Range("Q6").Value = 381.7084206403649
$ws.Range("R6").Value = 3435.375785763284
$ws.Range("S6").Value = 0.09180665363509669
$ws.Range("T6").Value = 0.09180665363509664

# Row 7: FAPs -> sCs (Col18a1/Gpc4)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col18a1"
$ws.Range("C7").Value = "Gpc4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.28572166666667
$ws.Range("H7").Value = 33.85716499999999
$ws.Range("I7").Value = 0.1500421083721963
$ws.Range("J7").Value = 0.1500421083721963
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.713974666666666
$ws.Range("N7").Value = 17.141924
$ws.Range("O7").Value = 0.1033705684511578
$ws.Range("P7").Value = 0.1033705684511578
$ws.Range("Q7").Value = 64.48632769838443
$ws.Range("R7").Value = 580.3769492854599
$ws.Range("S7").Value = 0.01550993803404416
$ws.Range("T7").Value = 0.01550993803404415

# Row 8: M1 -> ECs (Col18a1/Gpc4)
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Col18a1"
$ws.Range("C8").Value = "Gpc4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.240457
$ws.Range("H8").Value = 0.721371
$ws.Range("I8").Value = 0.003196842551895872
$ws.Range("J8").Value = 0.003196842551895872
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 15.740393
$ws.Range("N8").Value = 47.221179
$ws.Range("O8").Value = 0.2847568403735705
$ws.Range("P8").Value = 0.2847568403735705
$ws.Range("Q8").Value = 3.784887679601
$ws.Range("R8").Value = 34.063989116409
$ws.Range("S8").Value = 0.0009103227842496505
$ws.Range("T8").Value = 0.0009103227842496504

# Row 9: M1 -> FAPs (Col18a1/Gpc4)
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Col18a1"
$ws.Range("C9").Value = "Gpc4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.240457
$ws.Range("H9").Value = 0.721371
$ws.Range("I9").Value = 0.003196842551895872
$ws.Range("J9").Value = 0.003196842551895872
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 33.82224299999999
$ws.Range("N9").Value = 101.466729
$ws.Range("O9").Value = 0.6118725911752718
$ws.Range("P9").Value = 0.6118725911752717
$ws.Range("Q9").Value = 8.132795085050999
$ws.Range("R9").Value = 73.19515576545899
$ws.Range("S9").Value = 0.001956060335807896
$ws.Range("T9").Value = 0.001956060335807895

# Row 10: M1 -> sCs (Col18a1/Gpc4)
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Col18a1"
$ws.Range("C10").Value = "Gpc4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.240457
$ws.Range("H10").Value = 0.721371
$ws.Range("I10").Value = 0.003196842551895872
$ws.Range("J10").Value = 0.003196842551895872
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.713974666666666
$ws.Range("N10").Value = 17.141924
$ws.Range("O10").Value = 0.1033705684511578
$ws.Range("P10").Value = 0.1033705684511578
$ws.Range("Q10").Value = 1.373965206422667
$ws.Range("R10").Value = 12.365686857804
$ws.Range("S10").Value = 0.0003304594318383263
$ws.Range("T10").Value = 0.0003304594318383263

# Row 11: M2 -> ECs (Col18a1/Gpc4)
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col18a1"
$ws.Range("C11").Value = "Gpc4"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.191908
$ws.Range("H11").Value = 0.575724
$ws.Range("I11").Value = 0.002551390312817814
$ws.Range("J11").Value = 0.002551390312817813
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 15.740393
$ws.Range("N11").Value = 47.221179
$ws.Range("O11").Value = 0.2847568403735705
$ws.Range("P11").Value = 0.2847568403735705
$ws.Range("Q11").Value = 3.020707339844
$ws.Range("R11").Value = 27.186366058596
$ws.Range("S11").Value = 0.0007265258440377363
$ws.Range("T11").Value = 0.0007265258440377361

# Row 12: M2 -> FAPs (Col18a1/Gpc4)
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col18a1"
$ws.Range("C12").Value = "Gpc4"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.191908
$ws.Range("H12").Value = 0.575724
$ws.Range("I12").Value = 0.002551390312817814
$ws.Range("J12").Value = 0.002551390312817813
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 33.82224299999999
$ws.Range("N12").Value = 101.466729
$ws.Range("O12").Value = 0.6118725911752718
$ws.Range("P12").Value = 0.6118725911752717
$ws.Range("Q12").Value = 6.490759009643998
$ws.Range("R12").Value = 58.41683108679599
$ws.Range("S12").Value = 0.001561125801803323
$ws.Range("T12").Value = 0.001561125801803323

# Row 13: M2 -> sCs (Col18a1/Gpc4)
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col18a1"
$ws.Range("C13").Value = "Gpc4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.191908
$ws.Range("H13").Value = 0.575724
$ws.Range("I13").Value = 0.002551390312817814
$ws.Range("J13").Value = 0.002551390312817813
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 5.713974666666666
$ws.Range("N13").Value = 17.141924
$ws.Range("O13").Value = 0.1033705684511578
$ws.Range("P13").Value = 0.1033705684511578
$ws.Range("Q13").Value = 1.096557450330667
$ws.Range("R13").Value = 9.869017052976
$ws.Range("S13").Value = 0.0002637386669767548
$ws.Range("T13").Value = 0.0002637386669767548

# Row 14: sCs -> ECs (Col18a1/Gpc4)
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col18a1"
$ws.Range("C14").Value = "Gpc4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 16.24228466666667
$ws.Range("H14").Value = 48.726854
$ws.Range("I14").Value = 0.2159389278016688
$ws.Range("J14").Value = 0.2159389278016687
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 15.740393
$ws.Range("N14").Value = 47.221179
$ws.Range("O14").Value = 0.2847568403735705
$ws.Range("P14").Value = 0.2847568403735705
$ws.Range("Q14").Value = 255.6599438712073
$ws.Range("R14").Value = 2300.939494840866
$ws.Range("S14").Value = 0.06149008679445975
$ws.Range("T14").Value = 0.06149008679445973

# Row 15: sCs -> FAPs (Col18a1/Gpc4)
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col18a1"
$ws.Range("C15").Value = "Gpc4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 16.24228466666667
$ws.Range("H15").Value = 48.726854
$ws.Range("I15").Value = 0.2159389278016688
$ws.Range("J15").Value = 0.2159389278016687
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 33.82224299999999
$ws.Range("N15").Value = 101.466729
$ws.Range("O15").Value = 0.6118725911752718
$ws.Range("P15").Value = 0.6118725911752717
$ws.Range("Q15").Value = 549.3504988711738
$ws.Range("R15").Value = 4944.154489840565
$ws.Range("S15").Value = 0.132127111289617
$ws.Range("T15").Value = 0.132127111289617

# Row 16: sCs -> sCs (Col18a1/Gpc4)
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col18a1"
$ws.Range("C16").Value = "Gpc4"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 16.24228466666667
$ws.Range("H16").Value = 48.726854
$ws.Range("I16").Value = 0.2159389278016688
$ws.Range("J16").Value = 0.2159389278016687
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 5.713974666666666
$ws.Range("N16").Value = 17.141924
$ws.Range("O16").Value = 0.1033705684511578
$ws.Range("P16").Value = 0.1033705684511578
$ws.Range("Q16").Value = 92.80800311412177
$ws.Range("R16").Value = 835.2720280270959
$ws.Range("S16").Value = 0.02232172971759203
$ws.Range("T16").Value = 0.02232172971759202
